$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Save a template cell with the bold/border/centered style (style index 1) before clearing
$ws.Range("A3").Copy()
$ws.Range("Z1").PasteSpecial(-4122)

# Clear the whole data area (values + formatting)
$ws.Range("A1:B60").Clear()

# Restore the style template into Z2 as a stable copy source, then clean up Z1
$ws.Range("Z1").Copy()
$ws.Range("Z2").PasteSpecial(-4122)
$ws.Range("Z1").Clear()

function Set-StyledCell($addr, $val) {
    $ws.Range("Z2").Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value2 = $val
}

# Row 2
Set-StyledCell "B2" "PROJECT_#_0"

# Row 3
Set-StyledCell "A3" "E#3"
$ws.Range("B3").Value2 = 1

# Row 4
Set-StyledCell "A4" "E#14"
$ws.Range("B4").Value2 = 1

# Row 5
Set-StyledCell "A5" "C#0"
$ws.Range("B5").Value2 = 8.56

# Row 6
Set-StyledCell "A6" "C#1"
$ws.Range("B6").Value2 = 12.99

# Row 7
Set-StyledCell "A7" "C#2"
$ws.Range("B7").Value2 = 10.99

# Row 8
Set-StyledCell "A8" "C#3"
$ws.Range("B8").Value2 = 10.52

# Row 9
Set-StyledCell "A9" "FITNESS"
$ws.Range("B9").Value2 = 0.1779928703891769

# Row 13
Set-StyledCell "B13" "PROJECT_#_1"

# Row 14
Set-StyledCell "A14" "E#1"
$ws.Range("B14").Value2 = 1

# Row 15
Set-StyledCell "A15" "E#6"
$ws.Range("B15").Value2 = 1

# Row 16
Set-StyledCell "A16" "C#0"
$ws.Range("B16").Value2 = 8.5

# Row 17
Set-StyledCell "A17" "C#1"
$ws.Range("B17").Value2 = 13.62

# Row 18
Set-StyledCell "A18" "C#2"
$ws.Range("B18").Value2 = 8.790000000000001

# Row 19
Set-StyledCell "A19" "C#3"
$ws.Range("B19").Value2 = 7.23

# Row 20
Set-StyledCell "A20" "FITNESS"
$ws.Range("B20").Value2 = 0.1281089352209926

# Row 24
Set-StyledCell "B24" "PROJECT_#_2"

# Row 25
Set-StyledCell "A25" "E#5"
$ws.Range("B25").Value2 = 1

# Row 26
Set-StyledCell "A26" "E#7"
$ws.Range("B26").Value2 = 1

# Row 27
Set-StyledCell "A27" "E#9"
$ws.Range("B27").Value2 = 1

# Row 28
Set-StyledCell "A28" "C#0"
$ws.Range("B28").Value2 = 11.13

# Row 29
Set-StyledCell "A29" "C#1"
$ws.Range("B29").Value2 = 11.76

# Row 30
Set-StyledCell "A30" "C#2"
$ws.Range("B30").Value2 = 13.75

# Row 31
Set-StyledCell "A31" "C#3"
$ws.Range("B31").Value2 = 7.57

# Row 32
Set-StyledCell "A32" "FITNESS"
$ws.Range("B32").Value2 = 0.158706368553227

# Row 36
Set-StyledCell "B36" "PROJECT_#_3"

# Row 37
Set-StyledCell "A37" "E#11"
$ws.Range("B37").Value2 = 1

# Row 38
Set-StyledCell "A38" "E#12"
$ws.Range("B38").Value2 = 1

# Row 39
Set-StyledCell "A39" "C#0"
$ws.Range("B39").Value2 = 13.9

# Row 40
Set-StyledCell "A40" "C#1"
$ws.Range("B40").Value2 = 7.29

# Row 41
Set-StyledCell "A41" "C#2"
$ws.Range("B41").Value2 = 8.05

# Row 42
Set-StyledCell "A42" "C#3"
$ws.Range("B42").Value2 = 10.4

# Row 43
Set-StyledCell "A43" "FITNESS"
$ws.Range("B43").Value2 = 0.1272089621524845

# Row 47
Set-StyledCell "B47" "PROJECT_#_4"

# Row 48
Set-StyledCell "A48" "E#2"
$ws.Range("B48").Value2 = 1

# Row 49
Set-StyledCell "A49" "E#8"
$ws.Range("B49").Value2 = 1

# Row 50
Set-StyledCell "A50" "C#0"
$ws.Range("B50").Value2 = 10.29

# Row 51
Set-StyledCell "A51" "C#1"
$ws.Range("B51").Value2 = 8.14

# Row 52
Set-StyledCell "A52" "C#2"
$ws.Range("B52").Value2 = 3.44

# Row 53
Set-StyledCell "A53" "C#3"
$ws.Range("B53").Value2 = 11.15

# Row 54
Set-StyledCell "A54" "FITNESS"
$ws.Range("B54").Value2 = 0.2130081278849053

# Clean up helper/template cells
$ws.Range("Z2").Clear()

